$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the row above (row 4) into the new row 5 first so the
# new cells (in particular the date-formatted G5) reuse the existing style
# instead of minting a brand new one.
$ws.Range("A4:H4").Copy($ws.Range("A5:H5"))

$ws.Range("A5").Value = 9928.61
$ws.Range("B5").Value = 10044.120000000001
$ws.Range("C5").Value = 19.170000000000002
$ws.Range("D5").Value = 19.39
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = 1.1499999999999999
$ws.Range("G5").Value = 42609.505833333336
$ws.Range("H5").Value = $false
